# Split the final run of the paragraph
#   ", sem considerar consequências futuras."
# into two runs:
#   ", sem considerar consequências "  (unchanged formatting/run)
#   "futuras. "                        (new run, trailing space added)
# on slide 3, shape "Text Placeholder 2".

$p = $ppt.ActivePresentation

$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*consequ*ncias futuras.*") {
                $targetShape = $shape
                $targetSlide = $slide
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text

$searchWord = "futuras."
$pos = $fullText.IndexOf($searchWord)

# 1-based character index where "futuras." begins
$splitStart = $pos + 1
$splitLength = $searchWord.Length

$secondPart = $tr.Characters($splitStart, $splitLength)
$secondPart.Text = "futuras. "
